$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits right after "= K * n").
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Add the new sentence to the final (currently empty) numbered list item,
#    matching the Times New Roman formatting used throughout the document.
#    A throw-away placeholder is typed after the sentence so the new
#    "_GoBack" bookmark can be anchored *between* the sentence and the
#    placeholder (collapsed bookmarks placed at the very end of a
#    paragraph's text otherwise end up wrapping the preceding run instead
#    of following it). The placeholder is removed afterwards, leaving the
#    bookmark correctly positioned right after the sentence.
$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range
$sentence = "We were not asked to implement a locality version of Insertion Sort because it is already locality aware."
$placeholder = "ZZPLACEHOLDERZZ"
$r.InsertAfter($sentence + $placeholder)

$fullRange = $p.Range
$fullRange.Font.Name = "Times New Roman"
$fullRange.Font.NameBi = "Times New Roman"

$sentenceRange = $d.Content
$sentenceRange.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sentenceRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $sentenceRange) | Out-Null

$placeholderRange = $d.Content
$placeholderRange.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$placeholderRange.Delete()
